$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 29; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.NumberFormat = "@"
    $cell.Value = "10.01.20"
    $cell.Style = "Normal"
}

$ws.Range("H35").Select()
